# Generate Report for Handoff
# For the rows that are "Ready for handoff" (rows 4-7) in each localized-language
# sheet, bump the Priority to "ht" and refresh the Latest Handoff Datetime.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-25 20:29:29"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-25 20:29:34"
}
